$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.662.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.91%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.962.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.58%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'244.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.50%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +3.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'60.62"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +7.77%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.39%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0795"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.26%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'14.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +7.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.842"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'21.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.84%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.252.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "'  +3.97%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.961.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.71%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'36.558.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'70.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.30%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0856"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'230.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +3.86%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +6.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +4.05%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +10.60%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'160.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.32%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.45%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +12.17%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.69%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +6.36%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.18%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +21.66%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.90%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D39").Value = "'5.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +1.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +3.92%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'16.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.82%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.367.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.58%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.93%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'88.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.09%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.76%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'44.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.38%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +5.86%  "
$ws.Range("E51").Style = "Normal"
